# "cambiar tipo de datos" - change the number format ("data type") applied
# to each column of the table in Sheet1 (A1:E7):
#   - Column A (DIA), C (ANIO), D and E (numeric values) -> Number format "0"
#   - Column B (MES) and the header row (A1:E1, all text labels) -> Text format "@"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric columns: DIA, ANIO, and the two numeric data columns
$ws.Range("A2:A7").NumberFormat = "0"
$ws.Range("C2:C7").NumberFormat = "0"
$ws.Range("D2:D7").NumberFormat = "0"
$ws.Range("E2:E7").NumberFormat = "0"

# Text columns: the month column and the header row
$ws.Range("B1:B7").NumberFormat = "@"
$ws.Range("A1:E1").NumberFormat = "@"

# Move the active selection from G5 to G4
[void]$ws.Range("G4").Select()
